# formatting add leavequota Leavetypes to luis
$wb = $excel.ActiveWorkbook

$wsApp = $wb.Worksheets.Item("LeaveApplication")
$wsQuota = $wb.Worksheets.Item("LeaveQuota")

# --- LeaveApplication sheet ---
# Fix "Shortlist" flag for Family Care Leave (row 5): N -> Y
$wsApp.Range("C5").Value = "Y"

# Autofit / widen columns A:D to fit their content
$wsApp.Range("A1").ColumnWidth = 23.666666666666668
$wsApp.Range("B1").ColumnWidth = 18.333333333333332
$wsApp.Range("C1").ColumnWidth = 7.666666666666667
$wsApp.Range("D1").ColumnWidth = 9.5

# Update the view: scroll so row 7 is at the top, select A24
$wsApp.Activate()
$wsApp.Application.ActiveWindow.ScrollRow = 7
$wsApp.Range("A24").Select() | Out-Null

# --- LeaveQuota sheet ---
# Correct duplicate leave-type label: "Ext Maternity (UP-Non- SC)" -> "Ext Maternity (UP-Non SC)"
$wsQuota.Range("A13").Value = "Ext Maternity (UP-Non SC)"

# Re-activate LeaveQuota sheet (keeps it the selected tab) and update selection to C16
$wsQuota.Activate()
$wsQuota.Range("C16").Select() | Out-Null
